$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Meteo")

$ws.Range('E2').Value = '2026-02-06 08:47:57'
$ws.Range('H2').NumberFormat = "@"
$ws.Range('H2').Value = '95%'
$ws.Range('K2').Value = '0.2 MJ/m2'
$ws.Range('E3').Value = '2026-02-06 08:47:59'
$ws.Range('K3').Value = '0.1 MJ/m2'
$ws.Range('E4').Value = '2026-02-06 08:48:02'
$ws.Range('J4').Value = '994.7 hPa'
$ws.Range('K4').Value = '0.4 MJ/m2'
$ws.Range('O4').Value = '11.5 °C'
$ws.Range('E5').Value = '2026-02-06 08:48:05'
$ws.Range('J5').Value = '995.0 hPa'
$ws.Range('K5').Value = '0.4 MJ/m2'
$ws.Range('O5').Value = '7.7 °C'
$ws.Range('E6').Value = '2026-02-06 08:48:07'
$ws.Range('J6').Value = '996.1 hPa'
$ws.Range('K6').Value = '0.5 MJ/m2'
$ws.Range('E7').Value = '2026-02-06 08:48:10'
$ws.Range('J7').Value = '995.9 hPa'
$ws.Range('K7').Value = '0.5 MJ/m2'
$ws.Range('E8').Value = '2026-02-06 08:48:13'
$ws.Range('H8').NumberFormat = "@"
$ws.Range('H8').Value = '93%'
$ws.Range('K8').Value = '0.5 MJ/m2'
$ws.Range('M8').Value = '9.7 °C 8:29 TU'
$ws.Range('O8').Value = '5.5 °C'
$ws.Range('E9').Value = '2026-02-06 08:48:15'
$ws.Range('E10').Value = '2026-02-06 08:48:18'
$ws.Range('M10').Value = '6.9 °C 8:29 TU'
$ws.Range('O10').Value = '4.8 °C'
$ws.Range('E11').Value = '2026-02-06 08:48:20'
$ws.Range('J11').Value = '997.2 hPa'
$ws.Range('K11').Value = '0.1 MJ/m2'
$ws.Range('O11').Value = '3.3 °C'
$ws.Range('E12').Value = '2026-02-06 08:48:23'
$ws.Range('H12').NumberFormat = "@"
$ws.Range('H12').Value = '62%'
$ws.Range('K12').Value = '0.4 MJ/m2'
$ws.Range('E13').Value = '2026-02-06 08:48:26'
$ws.Range('E14').Value = '2026-02-06 08:48:28'
$ws.Range('I14').Value = '0.6 mm'
$ws.Range('K14').Value = '0.1 MJ/m2'
$ws.Range('O14').Value = '-3.9 °C'
$ws.Range('E15').Value = '2026-02-06 08:48:31'
$ws.Range('J15').Value = '995.2 hPa'
$ws.Range('K15').Value = '0.5 MJ/m2'
$ws.Range('O15').Value = '5.5 °C'
$ws.Range('E16').Value = '2026-02-06 08:48:34'
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H16').Value = '93%'
$ws.Range('K16').Value = '0.2 MJ/m2'
$ws.Range('E17').Value = '2026-02-06 08:48:37'
$ws.Range('J17').Value = '998.0 hPa'
$ws.Range('K17').Value = '0.3 MJ/m2'
$ws.Range('E18').Value = '2026-02-06 08:48:39'
$ws.Range('K18').Value = '0.3 MJ/m2'
$ws.Range('E19').Value = '2026-02-06 08:48:42'
$ws.Range('H19').NumberFormat = "@"
$ws.Range('H19').Value = '94%'
$ws.Range('J19').Value = '998.2 hPa'
$ws.Range('K19').Value = '0.4 MJ/m2'
$ws.Range('M19').Value = '9.6 °C 8:28 TU'
$ws.Range('O19').Value = '6.7 °C'
$ws.Range('E20').Value = '2026-02-06 08:48:45'
$ws.Range('K20').Value = '0.7 MJ/m2'
$ws.Range('E21').Value = '2026-02-06 08:48:47'
$ws.Range('J21').Value = '996.1 hPa'
$ws.Range('K21').Value = '0.4 MJ/m2'
$ws.Range('O21').Value = '4.1 °C'
$ws.Range('E22').Value = '2026-02-06 08:48:50'
$ws.Range('H22').NumberFormat = "@"
$ws.Range('H22').Value = '87%'
$ws.Range('K22').Value = '0.5 MJ/m2'
$ws.Range('O22').Value = '6.9 °C'
$ws.Range('E23').Value = '2026-02-06 08:48:53'
$ws.Range('J23').Value = '995.2 hPa'
$ws.Range('K23').Value = '0.4 MJ/m2'
$ws.Range('E24').Value = '2026-02-06 08:48:56'
$ws.Range('J24').Value = '994.2 hPa'
$ws.Range('K24').Value = '0.4 MJ/m2'
$ws.Range('M24').Value = '13.0 °C 8:07 TU'
$ws.Range('O24').Value = '12.2 °C'
$ws.Range('E25').Value = '2026-02-06 08:48:58'
$ws.Range('J25').Value = '997.3 hPa'
$ws.Range('K25').Value = '0.2 MJ/m2'
$ws.Range('L25').Value = '11.5 km/h - 298º 8:14 TU'
$ws.Range('E26').Value = '2026-02-06 08:49:00'
$ws.Range('K26').Value = '0.1 MJ/m2'
$ws.Range('O26').Value = '-1.8 °C'
$ws.Range('E27').Value = '2026-02-06 08:49:03'
$ws.Range('J27').Value = '995.0 hPa'
$ws.Range('K27').Value = '0.4 MJ/m2'
$ws.Range('O27').Value = '7.0 °C'
$ws.Range('E28').Value = '2026-02-06 08:49:06'
$ws.Range('H28').NumberFormat = "@"
$ws.Range('H28').Value = '93%'
$ws.Range('J28').Value = '998.5 hPa'
$ws.Range('O28').Value = '1.5 °C'
$ws.Range('E29').Value = '2026-02-06 08:49:08'
$ws.Range('K29').Value = '0.4 MJ/m2'
$ws.Range('E30').Value = '2026-02-06 08:49:11'
$ws.Range('K30').Value = '0.6 MJ/m2'
$ws.Range('O30').Value = '-3.9 °C'
$ws.Range('E31').Value = '2026-02-06 08:49:14'
$ws.Range('J31').Value = '997.9 hPa'
$ws.Range('N31').Value = '3.9 °C 8:09 TU'
$ws.Range('E32').Value = '2026-02-06 08:49:16'
$ws.Range('H32').NumberFormat = "@"
$ws.Range('H32').Value = '50%'
$ws.Range('J32').Value = '996.4 hPa'
$ws.Range('K32').Value = '0.4 MJ/m2'
$ws.Range('O32').Value = '14.6 °C'
$ws.Range('E33').Value = '2026-02-06 08:49:19'
$ws.Range('O33').Value = '6.2 °C'
$ws.Range('E34').Value = '2026-02-06 08:49:22'
$ws.Range('H34').NumberFormat = "@"
$ws.Range('H34').Value = '87%'
$ws.Range('K34').Value = '0.5 MJ/m2'
$ws.Range('O34').Value = '5.9 °C'
$ws.Range('E35').Value = '2026-02-06 08:49:24'
$ws.Range('K35').Value = '0.2 MJ/m2'
$ws.Range('N35').Value = '-3.5 °C 8:18 TU'
$ws.Range('E36').Value = '2026-02-06 08:49:27'
$ws.Range('J36').Value = '997.9 hPa'
$ws.Range('K36').Value = '0.3 MJ/m2'
